$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: 10_ulysses_6.tsp
$ws.Range("B4").Value = 29.00100016593933
$ws.Range("D4").Value = "[5][10][7][9, 6][2, 3][8, 1, 4]"

# Row 5: 10_ulysses_9.tsp
$ws.Range("B5").Value = 29.00100016593933

# Row 6: 14_burma_3.tsp
$ws.Range("B6").Value = 29.01699995994568

# Row 7: 14_burma_6.tsp
$ws.Range("B7").Value = 29.000999927520752
$ws.Range("C7").Value = "56.52842 / 42.740624"
$ws.Range("D7").Value = "[7][2, 13][12, 14, 6][1, 11, 8][10, 9][4, 3, 5]"

# Row 8: 14_burma_9.tsp
$ws.Range("B8").Value = 29.000999927520752
$ws.Range("D8").Value = "[7][11, 9][8, 1][13, 14][12, 6][5][10][4, 3][2]"

# Row 9: 22_ulysses_3.tsp
$ws.Range("B9").Value = 29.000999927520752
$ws.Range("C9").Value = "611.11393 / 445.92591"
$ws.Range("D9").Value = "[7, 5, 11, 10, 6, 9, 12][8, 17, 3, 16, 15, 4, 2][13, 20, 1, 14, 19, 21, 22, 18]"

# Row 10: 22_ulysses_6.tsp
$ws.Range("B10").Value = 29.000999927520752
$ws.Range("C10").Value = "185.71962 / 150.89245"
$ws.Range("D10").Value = "[21, 19, 1][7, 6, 5, 20][9, 11][8, 3, 4, 22, 18][15, 10, 12, 14][16, 13, 17, 2]"

# Row 11: 22_ulysses_9.tsp
$ws.Range("B11").Value = 29.0479998588562
$ws.Range("D11").Value = "[7, 1][2, 21][22, 16, 17][5, 11][15, 14, 12][4, 18, 8][9, 3][10, 6][20, 19, 13]"

# Row 12: 26_eil_3.tsp
$ws.Range("B12").Value = 29.003999948501587
$ws.Range("C12").Value = "3300.0705 / 2937.2667"
$ws.Range("D12").Value = "[23, 26, 7, 8, 18, 24, 1, 3][2, 9, 21, 5, 4, 10, 11, 20][19, 6, 25, 14, 22, 13, 12, 17, 15, 16]"

# Row 13: 26_eil_6.tsp
$ws.Range("B13").Value = 29.000999927520752
$ws.Range("C13").Value = "1504.4338 / 1037.622"
$ws.Range("D13").Value = "[4, 17, 5, 10][13, 19, 15, 14][16, 20, 21, 9, 11][1, 22, 26, 2, 3][18, 12, 6, 25][8, 23, 7, 24]"

# Row 14: 26_eil_9.tsp
$ws.Range("B14").Value = 29.000999927520752
$ws.Range("C14").Value = "1100.5689 / 762.82048"
$ws.Range("D14").Value = "[4, 25][18, 13][8, 24, 26][3, 1, 2][23, 22, 7][10, 9, 15][12, 19, 17][6, 14, 5][21, 11, 16, 20]"
